$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells I1 ("I0") and J1 ("IF"), matching the bold/centered/
# bordered header style already used by B1:H1 (style index 1 / "s=1").
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Row -> (I value, J value) for rows 2..79
$data = @(
    @(2,5,6),
    @(3,4,5),
    @(4,2,3),
    @(5,6,7),
    @(6,8,8),
    @(7,1,1),
    @(8,1,2),
    @(9,1,3),
    @(10,6,7),
    @(11,9,9),
    @(12,6,6),
    @(13,8,9),
    @(14,9,9),
    @(15,8,8),
    @(16,7,7),
    @(17,8,8),
    @(18,9,9),
    @(19,7,7),
    @(20,8,8),
    @(21,7,7),
    @(22,8,8),
    @(23,7,7),
    @(24,8,8),
    @(25,6,6),
    @(26,7,7),
    @(27,7,7),
    @(28,8,8),
    @(29,8,8),
    @(30,7,7),
    @(31,8,8),
    @(32,6,6),
    @(33,8,8),
    @(34,7,7),
    @(35,6,6),
    @(36,7,7),
    @(37,8,8),
    @(38,8,8),
    @(39,7,8),
    @(40,8,8),
    @(41,7,7),
    @(42,8,8),
    @(43,8,8),
    @(44,7,7),
    @(45,5,6),
    @(46,6,7),
    @(47,7,7),
    @(48,7,7),
    @(49,8,8),
    @(50,7,7),
    @(51,7,7),
    @(52,8,8),
    @(53,7,7),
    @(54,8,8),
    @(55,8,8),
    @(56,6,7),
    @(57,6,6),
    @(58,7,7),
    @(59,8,8),
    @(60,7,7),
    @(61,5,5),
    @(62,7,7),
    @(63,7,7),
    @(64,6,6),
    @(65,6,7),
    @(66,8,8),
    @(67,7,7),
    @(68,8,8),
    @(69,8,8),
    @(70,6,7),
    @(71,8,8),
    @(72,6,6),
    @(73,8,8),
    @(74,6,7),
    @(75,5,5),
    @(76,5,5),
    @(77,7,7),
    @(78,6,6),
    @(79,3,3)
)

foreach ($item in $data) {
    $row = $item[0]
    $ws.Cells.Item($row, 9).Value = $item[1]
    $ws.Cells.Item($row, 10).Value = $item[2]
}
